# Atualização de bases das ligas, do dia: 03-04-2024 às 22:09
# Updates odds for a few already-recorded fixtures, inserts one "new"
# fixture in date order (pushing the two following rows down by one),
# and appends two brand-new fixtures at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) In-place odds corrections for existing rows 235, 237, 238
# ---------------------------------------------------------------------

$ws.Range("R235").Value = 1.775
$ws.Range("S235").Value = 2.1
$ws.Range("U235").Value = 1.975
$ws.Range("V235").Value = 1.875

$ws.Range("N237").Value = 1.8
$ws.Range("P237").Value = 4.5
$ws.Range("Q237").Value = -0.5
$ws.Range("R237").Value = 1.825
$ws.Range("S237").Value = 2.025
$ws.Range("U237").Value = 1.925
$ws.Range("V237").Value = 1.925

$ws.Range("O238").Value = 3.5
$ws.Range("Q238").Value = -0.5
$ws.Range("R238").Value = 1.85
$ws.Range("S238").Value = 2
$ws.Range("U238").Value = 1.85
$ws.Range("V238").Value = 2

# ---------------------------------------------------------------------
# 2) Insert a new fixture row at row 239 (existing rows 239 & 240 shift
#    down to 240 & 241 respectively, unchanged)
# ---------------------------------------------------------------------

$ws.Rows.Item(239).Insert()

# restore the row-index column (A) styling on the freshly inserted row
# (bold, centered/top aligned, thin box border - same "style 1" used by
# every other data row) and the date number format on column E, by
# cloning the formats from an untouched reference row instead of
# re-building them property-by-property (avoids leaving behind unused
# intermediate cellXfs entries in styles.xml)
$ws.Range("A236").Copy()
$ws.Range("A239").PasteSpecial(-4122)
$ws.Range("E236").Copy()
$ws.Range("E239").PasteSpecial(-4122)

$ws.Range("A239").Value = 237
$ws.Range("B239").Value = 6775587
$ws.Range("C239").Value = "Poland Ekstraklasa"
$ws.Range("D239").Value = "Poland Ekstraklasa"
$ws.Range("E239").Value = 45388.625
$ws.Range("F239").Value = "Slask Wroclaw"
$ws.Range("G239").Value = "Warta Poznan"
$ws.Range("K239").Value = 1.727
$ws.Range("L239").Value = 3.8
$ws.Range("M239").Value = 4.5
$ws.Range("N239").Value = 1.909
$ws.Range("O239").Value = 3.6
$ws.Range("P239").Value = 3.75
$ws.Range("Q239").Value = -0.5
$ws.Range("R239").Value = 1.975
$ws.Range("S239").Value = 1.875
$ws.Range("T239").Value = 2
$ws.Range("U239").Value = 1.925
$ws.Range("V239").Value = 1.925
$ws.Range("W239").Value = 0
$ws.Range("X239").Value = 0
$ws.Range("Y239").Value = 0
$ws.Range("Z239").Value = 0
$ws.Range("AA239").Value = 0

# ---------------------------------------------------------------------
# 3) Append two brand-new fixture rows at the end (242 & 243)
# ---------------------------------------------------------------------

$ws.Range("A236").Copy()
$ws.Range("A242").PasteSpecial(-4122)
$ws.Range("E236").Copy()
$ws.Range("E242").PasteSpecial(-4122)

$ws.Range("A242").Value = 240
$ws.Range("B242").Value = 6775584
$ws.Range("C242").Value = "Poland Ekstraklasa"
$ws.Range("D242").Value = "Poland Ekstraklasa"
$ws.Range("E242").Value = 45389.52083333334
$ws.Range("F242").Value = "Legia Warsaw"
$ws.Range("G242").Value = "Jagiellonia Bialystok"
$ws.Range("K242").Value = 2
$ws.Range("L242").Value = 3.5
$ws.Range("M242").Value = 3.6
$ws.Range("N242").Value = 1.85
$ws.Range("O242").Value = 3.6
$ws.Range("P242").Value = 4
$ws.Range("Q242").Value = -0.5
$ws.Range("R242").Value = 1.9
$ws.Range("S242").Value = 1.95
$ws.Range("T242").Value = 2.75
$ws.Range("U242").Value = 1.925
$ws.Range("V242").Value = 1.925
$ws.Range("W242").Value = 0
$ws.Range("X242").Value = 0
$ws.Range("Y242").Value = 0
$ws.Range("Z242").Value = 0
$ws.Range("AA242").Value = 0

$ws.Range("A236").Copy()
$ws.Range("A243").PasteSpecial(-4122)
$ws.Range("E236").Copy()
$ws.Range("E243").PasteSpecial(-4122)

$ws.Range("A243").Value = 241
$ws.Range("B243").Value = 6775579
$ws.Range("C243").Value = "Poland Ekstraklasa"
$ws.Range("D243").Value = "Poland Ekstraklasa"
$ws.Range("E243").Value = 45390.58333333334
$ws.Range("F243").Value = "Zaglebie Lubin"
$ws.Range("G243").Value = "Gornik Zabrze"
$ws.Range("K243").Value = 2.2
$ws.Range("L243").Value = 3.4
$ws.Range("M243").Value = 3.2
$ws.Range("N243").Value = 2.2
$ws.Range("O243").Value = 3.4
$ws.Range("P243").Value = 3.2
$ws.Range("Q243").Value = -0.25
$ws.Range("R243").Value = 1.925
$ws.Range("S243").Value = 1.925
$ws.Range("T243").Value = 2.5
$ws.Range("U243").Value = 1.95
$ws.Range("V243").Value = 1.9
$ws.Range("W243").Value = 0
$ws.Range("X243").Value = 0
$ws.Range("Y243").Value = 0
$ws.Range("Z243").Value = 0
$ws.Range("AA243").Value = 0
